$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.957.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.740.31"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.68%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.91%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -1.81%  "
$ws.Range("E9").Value = "  -1.39%  "
$ws.Range("E10").Value = "  +4.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.75"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.83%  "
$ws.Range("E12").Value = "  -0.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.224.70"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.791.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.04%  "
$ws.Range("E16").Value = "  -1.50%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.742.77"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("E19").Value = "  -1.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "354.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.79%  "
$ws.Range("E21").Value = "  -3.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("E23").Value = "  -4.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.37"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.17%  "
$ws.Range("E25").Value = "  -0.70%  "
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("E27").Value = "  -1.25%  "
$ws.Range("E28").Value = "  -1.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("E30").Value = "  +2.93%  "
$ws.Range("E31").Value = "  +7.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "163.96"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.16%  "
$ws.Range("E33").Value = "  -1.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.66%  "
$ws.Range("E35").Value = "  +1.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("E37").Value = "  -0.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.987"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "345.60"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.34"
$ws.Range("D40").Style = "Normal"
$ws.Range("E41").Value = "  -1.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.56"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.85"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.50%  "
$ws.Range("E45").Value = "  -2.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "134.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.59%  "
$ws.Range("E47").Value = "  -1.90%  "
$ws.Range("E48").Value = "  -2.52%  "
$ws.Range("E49").Value = "  -1.44%  "
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.136.16"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.85%  "
